# First experiments completed, added results.
# Appends a new block of (Age, Prediction) observations below the existing
# table, letting column C's "A-B" margin-of-error formula fill in per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 is intentionally left blank (gap between the two data blocks),
# new data resumes at row 24 and runs through row 47.
$newData = @(
    @(24, 58, 44),
    @(25, 58, 50),
    @(26, 58, 54),
    @(27, 58, 48),
    @(28, 58, 46),
    @(29, 58, 51),
    @(30, 58, 47),
    @(31, 58, 58),
    @(32, 58, 50),
    @(33, 58, 45),
    @(34, 58, 45),
    @(35, 58, 50),
    @(36, 58, 42),
    @(37, 58, 44),
    @(38, 58, 43),
    @(39, 58, 50),
    @(40, 58, 50),
    @(41, 58, 51),
    @(42, 58, 55),
    @(43, 58, 62),
    @(44, 58, 51),
    @(45, 58, 52),
    @(46, 58, 42),
    @(47, 58, 49)
)

foreach ($entry in $newData) {
    $row = $entry[0]
    $age = $entry[1]
    $prediction = $entry[2]

    $ws.Cells.Item($row, 1).Value = $age
    $ws.Cells.Item($row, 2).Value = $prediction
    $ws.Cells.Item($row, 3).Formula = "=A$row-B$row"
}

# Leave the selection where the author left off working.
$ws.Range("D18").Select() | Out-Null
